$wb = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsEs = $wb.Worksheets.Item("es")

# --- Insert 2 new rows at row 7 in both sheets for newGame / continue ---
$wsEn.Range("A7:A8").EntireRow.Insert()
$wsEn.Cells.Item(7,1).Value = "newGame"
$wsEn.Cells.Item(7,2).Value = "NEW GAME"
$wsEn.Cells.Item(7,2).WrapText = $true
$wsEn.Cells.Item(8,1).Value = "continue"
$wsEn.Cells.Item(8,2).Value = "CONTINUE"
$wsEn.Cells.Item(8,2).WrapText = $true

$wsEs.Range("A7:A8").EntireRow.Insert()
$wsEs.Cells.Item(7,1).Value = "newGame"
$wsEs.Cells.Item(7,2).Value = "NUEVO JUEGO"
$wsEs.Cells.Item(8,1).Value = "continue"
$wsEs.Cells.Item(8,2).Value = "CONTINUAR"

# --- Append 2 new rows at the end (122, 123) for retry_desc / restart ---
$wsEn.Cells.Item(122,1).Value = "retry_desc"
$wsEn.Cells.Item(123,1).Value = "restart"
$wsEn.Cells.Item(122,2).Value = "It seems you had some trouble with this level. Press RESTART if you want to try again, or CONTINUE to go to the next lesson."
$wsEn.Cells.Item(123,2).Value = "RESTART"

$wsEs.Cells.Item(122,1).Value = "retry_desc"
$wsEs.Cells.Item(123,1).Value = "restart"
